# Apply "Added 'is_scaler' boolean to file specs" edit:
#  - tables sheet: new column E "is_scalar" (boolean) next to the existing
#    table/entity/snake_table/type columns.
#  - fields sheet: fix two VARCHAR -> DATETIME type mistakes (rows 119 & 124)
#    and drop the stale AutoFilter.
#  - active-tab/selection bookkeeping moves from the "fields" sheet back to
#    the "tables" sheet.

$wb = $excel.ActiveWorkbook
$wsTables = $wb.Worksheets.Item("tables")
$wsFields = $wb.Worksheets.Item("fields")

# --- tables sheet: add "is_scalar" boolean column -------------------------
$wsTables.Cells.Item(1, 5).Value = "is_scalar"
$wsTables.Cells.Item(1, 5).Font.Bold = $true

# Row -> is_scalar value for the 29 table specs currently listed (rows 2-30).
$scalarRows = @(13, 15)

for ($r = 2; $r -le 30; $r++) {
    if ($scalarRows -contains $r) {
        $wsTables.Cells.Item($r, 5).Value = $true
    } else {
        $wsTables.Cells.Item($r, 5).Value = $false
    }
}

# --- fields sheet: correct two mistaken VARCHAR types to DATETIME ---------
$wsFields.Cells.Item(119, 5).Value = "DATETIME"
$wsFields.Cells.Item(124, 5).Value = "DATETIME"

# Drop the leftover AutoFilter on the fields sheet.
$wsFields.AutoFilterMode = $false

# --- view bookkeeping: "tables" becomes the active/selected sheet ---------
$wsFields.Activate()
$wsFields.Range("A1").Select() | Out-Null

$wsTables.Activate()
$wsTables.Range("A1").Select() | Out-Null
